$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Notes" -> "Comments"
$ws.Range("H1").Value = "Comments"

# TC03 (News Article Navigation) block: fix step text, and shift the
# "categorized appropriately" expected-result text down from E8 to E9,
# clearing E8 and dropping the old E9 text ("...tagged with relevant keywords").
$ws.Range("D8").Value = "2. Locate the News section, click on it."
$ws.Range("E8").Value = ""
$ws.Range("E9").Value = " The news articles should be categorized appropriately (e.g., Ball-tampering, Technology in cricket, The future of ODI's)"

# TC05 (Match Schedules) block: wording tweak + typo fix.
$ws.Range("D14").Value = "2. Locate Live Scores menu option, under that locate Schedule option, click on it and check the schedule for a specific series or tournament"
$ws.Range("E15").Value = "Upon selecting appropriate filter conditions the corresponding data should be displayed."

# TC08 (Search Functionality) block: richer expected-result text.
$ws.Range("E23").Value = "The search results should display relevant information related to the player being searched.`nE.g. Player stats, Videos/Photos of the player etc."

# Remove the TC11 (Verify Advertisement Display) block entirely (rows 29-30).
$ws.Rows("29:30").Delete()

$ws.Range("D6").Select()
